# MMF_Tracker_2026_For_Website.xlsx - "Add files via upload" update
#
# The January sheet's Interest Earned (D) / Days Active (E) / Closing
# Balance (F) columns were refreshed with a new day's figures (one more
# day of accrued interest per MMF), and the Interest Earned / Closing
# Balance columns were given Excel's built-in "Comma" number style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January")

# --- Updated Days Active / Interest Earned / Closing Balance figures ---

# Row 4 - Tokyo (no activity this month)
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 5 - Berlin
$ws.Range("D5").Value = 19.032011880721566
$ws.Range("E5").Value = 26
$ws.Range("F5").Value = 6019.0320118807213

# Row 6 - Gandia
$ws.Range("D6").Value = 14.623502519883198
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 6014.623502519883

# Row 7 - Denver
$ws.Range("D7").Value = 17.2628332879999
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 9017.2628332879995

# Row 8 - Helsinki (no activity this month)
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0

# Row 9 - Oslo
$ws.Range("D9").Value = 1.9249393249216231
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 4001.9249393249215

# Row 10 - Moscow (no activity this month)
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# Row 11 - California (no activity this month)
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

# Row 12 - Lisbon
$ws.Range("D12").Value = 3.8667159572110834
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 4003.866715957211

# Row 13 - Stockholm (no activity this month)
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0

# Row 14 - Paris
$ws.Range("D14").Value = 0.3687639872599035
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 510.36876398725991

# Row 15 - Marseille
$ws.Range("D15").Value = 2.2905403364476014
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 2502.2905403364475

# Row 16 - Rio
$ws.Range("D16").Value = 16.071381864107565
$ws.Range("E16").Value = 19
$ws.Range("F16").Value = 6016.0713818641079

# Row 17 - Nairobi
$ws.Range("D17").Value = 24.103960362852376
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 20024.103960362852

# Row 18 - Tel Aviv
$ws.Range("D18").Value = 32.638870974089521
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 14032.63887097409

# Row 19 - Valencia
$ws.Range("D19").Value = 87.12738884250814
$ws.Range("E19").Value = 12
$ws.Range("F19").Value = 60087.127388842506

# Row 20 - Manilla
$ws.Range("D20").Value = 12.694495504278638
$ws.Range("E20").Value = 23
$ws.Range("F20").Value = 4012.6944955042786

# Row 21 - Scofield
$ws.Range("D21").Value = 2.4061741561520291
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 5002.4061741561518

# Row 22 - Bogota (no activity this month)
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0

# Row 23 - Emirates
$ws.Range("D23").Value = 14.829928412066284
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 7014.8299284120667

# Row 24 - Dublin (no activity this month)
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

# Row 25 - Doha
$ws.Range("D25").Value = 5.8084925895005446
$ws.Range("E25").Value = 12
$ws.Range("F25").Value = 4005.8084925895005

# Row 26 - TOTAL (Days Active column has no total; D/F are column sums)
$ws.Range("D26").Value = 255.04999999999998
$ws.Range("F26").Value = 152265.05000000002

# --- Apply the built-in "Comma" cell style to the Interest Earned and
#     Closing Balance columns (matches clicking the Comma Style button
#     on the Home ribbon for the D4:D26 and F4:F26 ranges) ---

$ws.Range("D4:D26").Style = "Comma"
$ws.Range("F4:F26").Style = "Comma"

# --- Restore the active selection left behind by the edit ---

$ws.Range("J31").Select() | Out-Null
